# Add a new "2021" data column (column R) to the table, mirroring the
# formatting of the preceding "2020" column (Q) cell-by-cell, and move
# the active selection from P9 down to P10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-YearColumn($row, $value) {
    $target = $ws.Range("R$row")
    $source = $ws.Range("Q$row")

    # Write the value first, then clone the source cell's formatting onto
    # it (copy / paste-special formats), matching how the column was
    # extended from the preceding year's cell.
    $target.Value = $value
    $source.Copy()
    $target.PasteSpecial(-4122)
}

Add-YearColumn 4 2021
Add-YearColumn 5 47.8
Add-YearColumn 6 20.7
Add-YearColumn 7 9.8000000000000007
Add-YearColumn 8 17.3

# Update the stored active cell / selection to match the authored diff.
$ws.Range("P10").Select()
